$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.568356333333332
$ws.Range("H2").Value = 28.705069
$ws.Range("I2").Value = 0.09973288675158326
$ws.Range("J2").Value = 0.09973288675158326
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.999936
$ws.Range("N2").Value = 41.999808
$ws.Range("O2").Value = 0.2017049292741484
$ws.Range("P2").Value = 0.2017049292741485
$ws.Range("Q2").Value = 133.9563762918613
$ws.Range("R2").Value = 1205.607386626752
$ws.Range("S2").Value = 0.02011661486853476
$ws.Range("T2").Value = 0.02011661486853476

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.568356333333332
$ws.Range("H3").Value = 28.705069
$ws.Range("I3").Value = 0.09973288675158326
$ws.Range("J3").Value = 0.09973288675158326
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("N3").Value = 113.242516
$ws.Range("O3").Value = 0.5438494785644407
$ws.Range("P3").Value = 0.5438494785644408
$ws.Range("Q3").Value = 361.1815817237338
$ws.Range("R3").Value = 3250.634235513604
$ws.Range("S3").Value = 0.05423967845557497
$ws.Range("T3").Value = 0.05423967845557498

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.568356333333332
$ws.Range("H4").Value = 28.705069
$ws.Range("I4").Value = 0.09973288675158326
$ws.Range("J4").Value = 0.09973288675158326
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 17.66056
$ws.Range("N4").Value = 52.98168
$ws.Range("O4").Value = 0.2544455921614109
$ws.Range("P4").Value = 0.2544455921614109
$ws.Range("Q4").Value = 168.9825311262133
$ws.Range("R4").Value = 1520.84278013592
$ws.Range("S4").Value = 0.02537659342747354
$ws.Range("T4").Value = 0.02537659342747354

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 75.94550566666668
$ws.Range("H5").Value = 227.836517
$ws.Range("I5").Value = 0.7915951551217724
$ws.Range("J5").Value = 0.7915951551217723
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.999936
$ws.Range("N5").Value = 41.999808
$ws.Range("O5").Value = 0.2017049292741484
$ws.Range("P5").Value = 0.2017049292741485
$ws.Range("Q5").Value = 1063.232218820971
$ws.Range("R5").Value = 9569.089969388737
$ws.Range("S5").Value = 0.1596686447775957
$ws.Range("T5").Value = 0.1596686447775957

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 75.94550566666668
$ws.Range("H6").Value = 227.836517
$ws.Range("I6").Value = 0.7915951551217724
$ws.Range("J6").Value = 0.7915951551217723
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 37.74750533333334
$ws.Range("N6").Value = 113.242516
$ws.Range("O6").Value = 0.5438494785644407
$ws.Range("P6").Value = 0.5438494785644408
$ws.Range("Q6").Value = 2866.753380195198
$ws.Range("R6").Value = 25800.78042175678
$ws.Range("S6").Value = 0.4305086123471134
$ws.Range("T6").Value = 0.4305086123471135

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 75.94550566666668
$ws.Range("H7").Value = 227.836517
$ws.Range("I7").Value = 0.7915951551217724
$ws.Range("J7").Value = 0.7915951551217723
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 17.66056
$ws.Range("N7").Value = 52.98168
$ws.Range("O7").Value = 0.2544455921614109
$ws.Range("P7").Value = 0.2544455921614109
$ws.Range("Q7").Value = 1341.240159556507
$ws.Range("R7").Value = 12071.16143600856
$ws.Range("S7").Value = 0.2014178979970633
$ws.Range("T7").Value = 0.2014178979970633

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.42596933333333
$ws.Range("H8").Value = 31.277908
$ws.Range("I8").Value = 0.1086719581266445
$ws.Range("J8").Value = 0.1086719581266445
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 13.999936
$ws.Range("N8").Value = 41.999808
$ws.Range("O8").Value = 0.2017049292741484
$ws.Range("P8").Value = 0.2017049292741485
$ws.Range("Q8").Value = 145.9629034046293
$ws.Range("R8").Value = 1313.666130641664
$ws.Range("S8").Value = 0.02191966962801804
$ws.Range("T8").Value = 0.02191966962801805

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.42596933333333
$ws.Range("H9").Value = 31.277908
$ws.Range("I9").Value = 0.1086719581266445
$ws.Range("J9").Value = 0.1086719581266445
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.74750533333334
$ws.Range("N9").Value = 113.242516
$ws.Range("O9").Value = 0.5438494785644407
$ws.Range("P9").Value = 0.5438494785644408
$ws.Range("Q9").Value = 393.5543330151698
$ws.Range("R9").Value = 3541.988997136528
$ws.Range("S9").Value = 0.05910118776175233
$ws.Range("T9").Value = 0.05910118776175233

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.42596933333333
$ws.Range("H10").Value = 31.277908
$ws.Range("I10").Value = 0.1086719581266445
$ws.Range("J10").Value = 0.1086719581266445
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.66056
$ws.Range("N10").Value = 52.98168
$ws.Range("O10").Value = 0.2544455921614109
$ws.Range("P10").Value = 0.2544455921614109
$ws.Range("Q10").Value = 184.1284569694933
$ws.Range("R10").Value = 1657.15611272544
$ws.Range("S10").Value = 0.0276511007368741
$ws.Range("T10").Value = 0.0276511007368741
